$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.505.01"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.39"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.11"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4857"
$ws.Range("E7").Value = "  +3.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2896"
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06710"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("B10").Value = "Litecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "109.03"
$ws.Range("E10").Value = "  +2.62%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.17"
$ws.Range("E11").Value = "  +6.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.914.78"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07552"
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.257"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6659"
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "270.99"
$ws.Range("E16").Value = "  -5.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.506.77"

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007528"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.169.68"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.491"
$ws.Range("E22").Value = "  +5.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9976"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.391"
$ws.Range("E24").Value = "  +3.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.401"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.68"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -4.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.093"
$ws.Range("E28").Value = "  +2.99%  "

$ws.Range("E29").Value = "  -2.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.116"
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.033"
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04980"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7271"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02028"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.667"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.93"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.011"
$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4419"
$ws.Range("E42").Value = "  +5.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8662"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.847"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.64"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.364"
$ws.Range("E47").Value = "  +3.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.273"
$ws.Range("E48").Value = "  +1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1242"
$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.26"
$ws.Range("E50").Value = "  -9.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.460"
$ws.Range("E51").Value = "  +6.28%  "
